$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "C3"
$ws.Range("C2").Value = "Itgam"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 4.153678333333333
$ws.Range("H2").Value = 12.461035
$ws.Range("I2").Value = 0.02267710693885585
$ws.Range("J2").Value = 0.02267710693885585
$ws.Range("K2").Value = 1.0
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.142723
$ws.Range("N2").Value = 0.428169
$ws.Range("O2").Value = 0.000470790275435748
$ws.Range("P2").Value = 0.0004707902754357479
$ws.Range("Q2").Value = 0.5928254327683333
$ws.Range("R2").Value = 5.335428894915
$ws.Range("S2").Value = 0.00001067616142182986
$ws.Range("T2").Value = 0.00001067616142182985

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "C3"
$ws.Range("C3").Value = "Itgam"
$ws.Range("D3").Value = "M1"
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 4.153678333333333
$ws.Range("H3").Value = 12.461035
$ws.Range("I3").Value = 0.02267710693885585
$ws.Range("J3").Value = 0.02267710693885585
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 148.0881626666667
$ws.Range("N3").Value = 444.264488
$ws.Range("O3").Value = 0.4884879584272602
$ws.Range("P3").Value = 0.4884879584272601
$ws.Range("Q3").Value = 615.1105926916755
$ws.Range("R3").Value = 5535.99533422508
$ws.Range("S3").Value = 0.01107749367159835
$ws.Range("T3").Value = 0.01107749367159835

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "C3"
$ws.Range("C4").Value = "Itgam"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 4.153678333333333
$ws.Range("H4").Value = 12.461035
$ws.Range("I4").Value = 0.02267710693885585
$ws.Range("J4").Value = 0.02267710693885585
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 154.9253336666667
$ws.Range("N4").Value = 464.776001
$ws.Range("O4").Value = 0.5110412512973043
$ws.Range("P4").Value = 0.5110412512973042
$ws.Range("Q4").Value = 643.5100017356705
$ws.Range("R4").Value = 5791.590015621034
$ws.Range("S4").Value = 0.01158893710583567
$ws.Range("T4").Value = 0.01158893710583567

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "C3"
$ws.Range("C5").Value = "Itgam"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 140.5890806666667
$ws.Range("H5").Value = 421.767242
$ws.Range("I5").Value = 0.7675494732291734
$ws.Range("J5").Value = 0.7675494732291734
$ws.Range("K5").Value = 1.0
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.142723
$ws.Range("N5").Value = 0.428169
$ws.Range("O5").Value = 0.000470790275435748
$ws.Range("P5").Value = 0.0004707902754357479
$ws.Range("Q5").Value = 20.06529535998867
$ws.Range("R5").Value = 180.587658239898
$ws.Range("S5").Value = 0.0003613548279121258
$ws.Range("T5").Value = 0.0003613548279121258

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "C3"
$ws.Range("C6").Value = "Itgam"
$ws.Range("D6").Value = "M1"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 140.5890806666667
$ws.Range("H6").Value = 421.767242
$ws.Range("I6").Value = 0.7675494732291734
$ws.Range("J6").Value = 0.7675494732291734
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 148.0881626666667
$ws.Range("N6").Value = 444.264488
$ws.Range("O6").Value = 0.4884879584272602
$ws.Range("P6").Value = 0.4884879584272601
$ws.Range("Q6").Value = 20819.57864692245
$ws.Range("R6").Value = 187376.2078223021
$ws.Range("S6").Value = 0.3749386751696379
$ws.Range("T6").Value = 0.3749386751696379

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "C3"
$ws.Range("C7").Value = "Itgam"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 140.5890806666667
$ws.Range("H7").Value = 421.767242
$ws.Range("I7").Value = 0.7675494732291734
$ws.Range("J7").Value = 0.7675494732291734
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 154.9253336666667
$ws.Range("N7").Value = 464.776001
$ws.Range("O7").Value = 0.5110412512973043
$ws.Range("P7").Value = 0.5110412512973042
$ws.Range("Q7").Value = 21780.81023217325
$ws.Range("R7").Value = 196027.2920895592
$ws.Range("S7").Value = 0.3922494432316235
$ws.Range("T7").Value = 0.3922494432316235

# Row 8
$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "C3"
$ws.Range("C8").Value = "Itgam"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 7.122973333333334
$ws.Range("H8").Value = 21.36892
$ws.Range("I8").Value = 0.03888804453304686
$ws.Range("J8").Value = 0.03888804453304686
$ws.Range("K8").Value = 1.0
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.142723
$ws.Range("N8").Value = 0.428169
$ws.Range("O8").Value = 0.000470790275435748
$ws.Range("P8").Value = 0.0004707902754357479
$ws.Range("Q8").Value = 1.016612123053334
$ws.Range("R8").Value = 9.149509107480002
$ws.Range("S8").Value = 0.00001830811319687077
$ws.Range("T8").Value = 0.00001830811319687076

# Row 9
$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "C3"
$ws.Range("C9").Value = "Itgam"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 7.122973333333334
$ws.Range("H9").Value = 21.36892
$ws.Range("I9").Value = 0.03888804453304686
$ws.Range("J9").Value = 0.03888804453304686
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 148.0881626666667
$ws.Range("N9").Value = 444.264488
$ws.Range("O9").Value = 0.4884879584272602
$ws.Range("P9").Value = 0.4884879584272601
$ws.Range("Q9").Value = 1054.828033656996
$ws.Range("R9").Value = 9493.452302912961
$ws.Range("S9").Value = 0.01899634148117644
$ws.Range("T9").Value = 0.01899634148117644

# Row 10
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "C3"
$ws.Range("C10").Value = "Itgam"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 7.122973333333334
$ws.Range("H10").Value = 21.36892
$ws.Range("I10").Value = 0.03888804453304686
$ws.Range("J10").Value = 0.03888804453304686
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 154.9253336666667
$ws.Range("N10").Value = 464.776001
$ws.Range("O10").Value = 0.5110412512973043
$ws.Range("P10").Value = 0.5110412512973042
$ws.Range("Q10").Value = 1103.529020365436
$ws.Range("R10").Value = 9931.76118328892
$ws.Range("S10").Value = 0.01987339493867356
$ws.Range("T10").Value = 0.01987339493867355

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "C3"
$ws.Range("C11").Value = "Itgam"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3.0
$ws.Range("F11").Value = 1.0
$ws.Range("G11").Value = 30.91341533333334
$ws.Range("H11").Value = 92.74024600000001
$ws.Range("I11").Value = 0.1687725358349285
$ws.Range("J11").Value = 0.1687725358349285
$ws.Range("K11").Value = 1.0
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.142723
$ws.Range("N11").Value = 0.428169
$ws.Range("O11").Value = 0.000470790275435748
$ws.Range("P11").Value = 0.0004707902754357479
$ws.Range("Q11").Value = 4.412055376619334
$ws.Range("R11").Value = 39.70849838957401
$ws.Range("S11").Value = 0.00007945646863171565
$ws.Range("T11").Value = 0.00007945646863171563

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "C3"
$ws.Range("C12").Value = "Itgam"
$ws.Range("D12").Value = "M1"
$ws.Range("E12").Value = 3.0
$ws.Range("F12").Value = 1.0
$ws.Range("G12").Value = 30.91341533333334
$ws.Range("H12").Value = 92.74024600000001
$ws.Range("I12").Value = 0.1687725358349285
$ws.Range("J12").Value = 0.1687725358349285
$ws.Range("K12").Value = 3.0
$ws.Range("L12").Value = 1.0
$ws.Range("M12").Value = 148.0881626666667
$ws.Range("N12").Value = 444.264488
$ws.Range("O12").Value = 0.4884879584272602
$ws.Range("P12").Value = 0.4884879584272601
$ws.Range("Q12").Value = 4577.910878464895
$ws.Range("R12").Value = 41201.19790618405
$ws.Range("S12").Value = 0.08244335146859584
$ws.Range("T12").Value = 0.08244335146859583

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "C3"
$ws.Range("C13").Value = "Itgam"
$ws.Range("D13").Value = "M2"
$ws.Range("E13").Value = 3.0
$ws.Range("F13").Value = 1.0
$ws.Range("G13").Value = 30.91341533333334
$ws.Range("H13").Value = 92.74024600000001
$ws.Range("I13").Value = 0.1687725358349285
$ws.Range("J13").Value = 0.1687725358349285
$ws.Range("K13").Value = 3.0
$ws.Range("L13").Value = 1.0
$ws.Range("M13").Value = 154.9253336666667
$ws.Range("N13").Value = 464.776001
$ws.Range("O13").Value = 0.5110412512973043
$ws.Range("P13").Value = 0.5110412512973042
$ws.Range("Q13").Value = 4789.271185292917
$ws.Range("R13").Value = 43103.44066763625
$ws.Range("S13").Value = 0.086249727897701
$ws.Range("T13").Value = 0.08624972789770098

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "C3"
$ws.Range("C14").Value = "Itgam"
$ws.Range("D14").Value = "FAPs"
$ws.Range("E14").Value = 3.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 0.3870006666666667
$ws.Range("H14").Value = 1.161002
$ws.Range("I14").Value = 0.002112839463995207
$ws.Range("J14").Value = 0.002112839463995207
$ws.Range("K14").Value = 1.0
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.142723
$ws.Range("N14").Value = 0.428169
$ws.Range("O14").Value = 0.000470790275435748
$ws.Range("P14").Value = 0.0004707902754357479
$ws.Range("Q14").Value = 0.05523389614866668
$ws.Range("R14").Value = 0.4971050653380001
$ws.Range("S14").Value = 0.000000994704273205822
$ws.Range("T14").Value = 0.0000009947042732058215

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "C3"
$ws.Range("C15").Value = "Itgam"
$ws.Range("D15").Value = "M1"
$ws.Range("E15").Value = 3.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 0.3870006666666667
$ws.Range("H15").Value = 1.161002
$ws.Range("I15").Value = 0.002112839463995207
$ws.Range("J15").Value = 0.002112839463995207
$ws.Range("K15").Value = 3.0
$ws.Range("L15").Value = 1.0
$ws.Range("M15").Value = 148.0881626666667
$ws.Range("N15").Value = 444.264488
$ws.Range("O15").Value = 0.4884879584272602
$ws.Range("P15").Value = 0.4884879584272601
$ws.Range("Q15").Value = 57.31021767744179
$ws.Range("R15").Value = 515.791959096976
$ws.Range("S15").Value = 0.001032096636251565
$ws.Range("T15").Value = 0.001032096636251565

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "C3"
$ws.Range("C16").Value = "Itgam"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 0.3870006666666667
$ws.Range("H16").Value = 1.161002
$ws.Range("I16").Value = 0.002112839463995207
$ws.Range("J16").Value = 0.002112839463995207
$ws.Range("K16").Value = 3.0
$ws.Range("L16").Value = 1.0
$ws.Range("M16").Value = 154.9253336666667
$ws.Range("N16").Value = 464.776001
$ws.Range("O16").Value = 0.5110412512973043
$ws.Range("P16").Value = 0.5110412512973042
$ws.Range("Q16").Value = 59.95620741255578
$ws.Range("R16").Value = 539.605866713002
$ws.Range("S16").Value = 0.001079748123470436
$ws.Range("T16").Value = 0.001079748123470436
